$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.520.38'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.826.53'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''316.56'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '''0.5175'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '''0.3888'
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("D9").Value = '''0.08407'
$ws.Range("E9").Value = '  +8.91%  '
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("D11").Value = '''41.93'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '''6.425'
$ws.Range("E12").Value = '  +2.25%  '
$ws.Range("D13").Value = '''21.16'
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '''7.521'
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").Value = '1.824.45'
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '''0.00001129'
$ws.Range("E17").Value = '  +4.52%  '
$ws.Range("D18").Value = '''93.72'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '''0.06619'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '''17.78'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '''1.002'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '''6.078'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = '28.560.88'
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("D25").Value = '''2.280'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").Value = '''21.15'
$ws.Range("E26").Value = '  +2.57%  '
$ws.Range("D27").Value = '''159.38'
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").Value = '2.033.53'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").Value = '''125.62'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").Value = '''0.1098'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").Value = '''1.100'
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").Value = '''5.739'
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").Value = '''0.07545'
$ws.Range("E34").Value = '  +4.41%  '
$ws.Range("D35").Value = '''3.666'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = '''0.02370'
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("D38").Value = '''5.232'
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").Value = '''8.781'
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("D40").Value = '''11.46'
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("D41").Value = '''0.6344'
$ws.Range("E41").Value = '  +1.52%  '
$ws.Range("D42").Value = '''1.193'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").Value = '''1.400'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").Value = '''13.56'
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.6007'
$ws.Range("E45").Value = '  +1.60%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.779'
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").Value = '''127.11'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("D48").Value = '''1.992'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("D50").Value = '''0.06983'
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").Value = '''74.62'
$ws.Range("E51").Value = '  +1.22%  '
